$wb = $excel.ActiveWorkbook

# "SE Results" is the template sheet for the decay-chain contamination
# calculations. Duplicate it to build a new scenario sheet for a 1 part-per-
# billion (1ppb) contamination level, placed after the last sheet.
$src = $wb.Worksheets.Item("SE Results")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy($null, $lastSheet)

$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "1ppb"

# Update the overall-contamination (per atom) inputs on the new sheet to
# 1e-9 (1ppb) for all three decay chains (232-Th, 238-U, 235-U blocks).
$ws.Range("D2").Value = 0.000000001
$ws.Range("D11").Value = 0.000000001
$ws.Range("D26").Value = 0.000000001

# Restore the selection on the original template sheet and leave the new
# "1ppb" sheet as the active tab with its own selection.
$src.Range("D28").Select()
$ws.Select()
$ws.Range("D27").Select()
